# Weekly update: prepend a new "Acelga" price record (most recent date)
# for Terminal La Palmera de La Serena, pushing the rest of the
# historical rows down by two rows (one for "Primera", one for
# "Segunda" quality).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the existing block (old row 415),
# shifting all data currently in rows 415:512 down to 417:514.
$ws.Rows("415:416").Insert()

# ---- New row 415 (Calidad = Primera) ----
$ws.Cells.Item(415, 1).Value  = 8
$ws.Cells.Item(415, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(415, 3).Value  = "Coquimbo"
$ws.Cells.Item(415, 4).Value  = 44855
$ws.Cells.Item(415, 5).Value  = 4
$ws.Cells.Item(415, 6).Value  = 100112009
$ws.Cells.Item(415, 7).Value  = "Acelga"
$ws.Cells.Item(415, 8).Value  = "Sin especificar"
$ws.Cells.Item(415, 9).Value  = "Primera"
$ws.Cells.Item(415, 10).Value = 2200
$ws.Cells.Item(415, 11).Value = 650
$ws.Cells.Item(415, 12).Value = 700
$ws.Cells.Item(415, 13).Value = 675
$ws.Cells.Item(415, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(415, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(415, 16).Value = 338
$ws.Cells.Item(415, 17).Value = 2
$ws.Cells.Item(415, 18).Value = "Hortaliza"

# ---- New row 416 (Calidad = Segunda) ----
$ws.Cells.Item(416, 1).Value  = 8
$ws.Cells.Item(416, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(416, 3).Value  = "Coquimbo"
$ws.Cells.Item(416, 4).Value  = 44855
$ws.Cells.Item(416, 5).Value  = 4
$ws.Cells.Item(416, 6).Value  = 100112009
$ws.Cells.Item(416, 7).Value  = "Acelga"
$ws.Cells.Item(416, 8).Value  = "Sin especificar"
$ws.Cells.Item(416, 9).Value  = "Segunda"
$ws.Cells.Item(416, 10).Value = 1500
$ws.Cells.Item(416, 11).Value = 550
$ws.Cells.Item(416, 12).Value = 600
$ws.Cells.Item(416, 13).Value = 575
$ws.Cells.Item(416, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(416, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(416, 16).Value = 288
$ws.Cells.Item(416, 17).Value = 2
$ws.Cells.Item(416, 18).Value = "Hortaliza"

# Make sure the date column keeps the same numeric (date) format as the
# rest of column D.
$ws.Range("D415:D416").NumberFormat = $ws.Range("D417").NumberFormat
